# "save up to 5 images URLs"
#
# The #property row (row 6, with #object_type_xsd in row 7 / #property_context
# in row 8) used column F for the single "image URL" property
# (http://schema.org/image). This edit makes room for up to 5 image URLs by
# inserting 4 more columns right after F, each carrying the same
# #property / #object_type_xsd / #property_context values as column F. The
# columns that used to sit at G:J (locatedIn / Assertion-marker / city / town
# / prefecture-ish metadata) shift right to K:N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new columns at G (pushes old G:J -> K:N, keeps all their values,
# formats and styles intact).
$ws.Range("G1:J1").EntireColumn.Insert()

# The new G:J columns should repeat column F's values for the property
# header rows (6 = property name/URI, 7 = data type, 8 = context).
$ws.Range("G6:J6").Value = $ws.Range("F6").Value()
$ws.Range("G7:J7").Value = $ws.Range("F7").Value()
$ws.Range("G8:J8").Value = $ws.Range("F8").Value()

# Keep the new columns' width in line with column F's (image URL) width.
$ws.Range("G1:J1").EntireColumn.ColumnWidth = $ws.Range("F1").EntireColumn.ColumnWidth
